$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "Reihenfolge" column values for the three data rows (H2:H4).
$ws.Range("H2").ClearContents()
$ws.Range("H3").ClearContents()
$ws.Range("H4").ClearContents()
